$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the numeric-looking "Price" column cells that change,
# so Excel COM does not coerce the assigned strings into floating-point numbers
# and we keep the exact textual representation (leading/trailing zeros, etc.)
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D40", "D41", "D42", "D44", "D45", "D47", "D48")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '243.84'
$ws.Range('D3').Value = '25.09'
$ws.Range('D4').Value = '5.188'
$ws.Range('D5').Value = '0.05741'
$ws.Range('D6').Value = '6.519'
$ws.Range('D11').Value = '0.06958'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.02828'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09373'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001516'
$ws.Range('E14').Value = '13BitForexTokenBF'
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').Value = '0.0006010'
$ws.Range('E15').Value = '14OneONE'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '0.006231'
$ws.Range('E16').Value = '15TigerCashTCH'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '3.504'
$ws.Range('E17').Value = '16LEOLEO'
$ws.Range('B18').Value = 'BTSEToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').Value = '2.092'
$ws.Range('E18').Value = '17BTSETokenBTSE'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').Value = '0.3185'
$ws.Range('E19').Value = '18BitpandaEcosystemTokenBEST'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').Value = '0.03127'
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('D22').Value = '3.749'
$ws.Range('D23').Value = '0.04673'
$ws.Range('D25').Value = '0.001233'
$ws.Range('D26').Value = '0.004265'
$ws.Range('D27').Value = '0.00008698'
$ws.Range('E28').Value = '27UpBotsUBXTWorstin24h'
$ws.Range('D40').Value = '0.03610'
$ws.Range('D41').Value = '0.006319'
$ws.Range('E41').Value = '40KickTokenKICKBestin24h'
$ws.Range('D42').Value = '0.1048'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('D44').Value = '0.007332'
$ws.Range('D45').Value = '0.00005298'
$ws.Range('D47').Value = '0.3440'
$ws.Range('D48').Value = '0.002280'
